# Auto update Excel log
# Appends new sensor/log rows to the Proximity, mmWave and Camera sheets.

function SetTextCell($ws, $r, $c, $text) {
    # Force text storage so date-looking strings (e.g. "2026-02-01") are not
    # auto-converted into Excel date serial numbers.
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function AppendRows($ws, $startRow, $rows) {
    $r = $startRow
    foreach ($row in $rows) {
        for ($c = 1; $c -le 6; $c++) {
            SetTextCell $ws $r $c $row[$c - 1]
        }
        $r++
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Proximity sheet: two new door ENTER/EXIT events (rows 49-50)
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "14:24:12", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "14:24:14", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)

AppendRows $wsProximity 49 $proximityRows

# ---------------------------------------------------------------------
# mmWave sheet: seven new presence-detection events (rows 7-13)
# ---------------------------------------------------------------------
$wsMmWave = $wb.Worksheets.Item("mmWave")

$mmWaveRows = @(
    @("2026-02-01", "14:24:13", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "14:24:15", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "14:24:25", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "14:24:36", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "14:24:46", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "14:24:57", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-02-01", "14:25:07", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

AppendRows $wsMmWave 7 $mmWaveRows

# ---------------------------------------------------------------------
# Camera sheet: one new "Image Captured" event (row 31)
# ---------------------------------------------------------------------
$wsCamera = $wb.Worksheets.Item("Camera")

# NOTE: the unary comma forces this to stay a one-element array of rows
# instead of being unwrapped into a single flat row by PowerShell.
$cameraRows = ,@("2026-02-01", "14:24:14", "14:00", "Living Room Main Door", "Image Captured", "Active")

AppendRows $wsCamera 31 $cameraRows
